$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "mod_Repairing" -> "mod_Repair" (set this first so it gets the lower
# shared-string index, matching the original authoring order)
$ws.Range("C22").Value = "mod_Repair"
$ws.Range("D22").Value = "mod_Repair"

# Rename "mod_Repowering" -> "mod_Reuse"
$ws.Range("C21").Value = "mod_Reuse"
$ws.Range("D21").Value = "mod_Reuse"

# Update the combined "variables" lists that referenced the old names
$ws.Range("D23").Value = "mod_EOL_collection_eff, mod_EOL_collected_recycled, mod_Repair, mod_Reuse"
$ws.Range("D24").Value = "mod_EOL_collection_eff, mod_EOL_collected_recycled, mod_Repair, mod_Reuse, mod_reliability_t50, mod_reliability_t90, mod_lifetime"

# Update the selection/active cell shown in the sheet view
$ws.Range("C22").Select()
